$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 detail cells: fill in the new template placeholder fields
# (pilihan_asal, keterangan_pengadaan, merk, nobukti, bahan, ukuran, satuan)
$ws.Range("H15").Value = "[a.pilihan_asal]"
$ws.Range("O15").Value = "[a.keterangan_pengadaan]"
$ws.Range("E15").Value = "[a.merk]"
$ws.Range("F15").Value = "[a.nobukti]"
$ws.Range("G15").Value = "[a.bahan]"
$ws.Range("J15").Value = "[a.ukuran]"
$ws.Range("K15").Value = "[a.satuan]"

# Update the view: scroll so column C is left-most and select O15
$ws.Range("O15").Select()
$excel.ActiveWindow.ScrollColumn = 3
